$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.281.21"
$ws.Range("E2").Value = "  +1.62%  "

$ws.Range("D3").Value = "2.003.15"
$ws.Range("E3").Value = "  +5.43%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.79"
$ws.Range("E5").Value = "  -1.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.657"
$ws.Range("E6").Value = "  -4.96%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.18"
$ws.Range("E8").Value = "  +1.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.79"
$ws.Range("E9").Value = "  +9.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.362"
$ws.Range("E10").Value = "  +0.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0710"
$ws.Range("E11").Value = "  -6.44%  "

$ws.Range("E12").Value = "  -0.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.32"
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").Value = "2.295.88"
$ws.Range("E14").Value = "  +5.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.801"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").Value = "2.003.09"
$ws.Range("E16").Value = "  +5.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.86"
$ws.Range("E17").Value = "  -3.54%  "

$ws.Range("D18").Value = "36.370.98"
$ws.Range("E18").Value = "  +1.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.00"
$ws.Range("E19").Value = "  -3.87%  "

$ws.Range("D20").Value = "0.0₃0811"
$ws.Range("E20").Value = "  -2.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.74"
$ws.Range("E21").Value = "  -2.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.25"
$ws.Range("E22").Value = "  -4.25%  "

$ws.Range("E23").Value = "  -6.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  -9.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.01"
$ws.Range("E26").Value = "  -1.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.56"
$ws.Range("E27").Value = "  -1.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.53"
$ws.Range("E28").Value = "  +6.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.92"
$ws.Range("E29").Value = "  -11.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.120"
$ws.Range("E30").Value = "  -6.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.90"
$ws.Range("E31").Value = "  +58.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.32"
$ws.Range("E32").Value = "  -1.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0580"
$ws.Range("E33").Value = "  -3.99%  "

$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("E35").Value = "  +0.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.95"
$ws.Range("E36").Value = "  -7.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0821"
$ws.Range("E37").Value = "  +11.02%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.846"
$ws.Range("E38").Value = "  -1.63%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.09"
$ws.Range("E39").Value = "  +6.02%  "

$ws.Range("E40").Value = "  -11.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0214"
$ws.Range("E41").Value = "  -5.48%  "

$ws.Range("E42").Value = "  +1.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "94.98"
$ws.Range("E43").Value = "  -4.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.76"
$ws.Range("E44").Value = "  +15.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.93"

$ws.Range("D46").Value = "1.306.09"
$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0812"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("E48").Value = "  +0.79%  "

$ws.Range("D49").Value = "2.189.08"
$ws.Range("E49").Value = "  +5.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.16"
$ws.Range("E50").Value = "  -8.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.82"
$ws.Range("E51").Value = "  +12.92%  "
